# Remove the literal "<exp>c</exp>" transcription markup that follows
# "&amp;" in the "graine de chou navette, &amp;<exp>c</exp>." paragraph.
# This deletes the <exp>, c, and </exp> runs entirely (their combined
# text), leaving "&amp;. Ce faict, ..." behind.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("<exp>c</exp>", $true, $false, $false, `
                                  $false, $false, $true, 1, $false, "", `
                                  2)
